# Update crypto price (D) and 1h-volume-change (E) columns per the latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.238.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "'1.674.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("E4").Value = "  -0.62%  "
$ws.Range("D5").Value = "'211.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.17%  "
$ws.Range("D6").Value = "'0.5268"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.82%  "
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("D8").Value = "'0.2653"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.00%  "
$ws.Range("D9").Value = "'0.06280"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.63%  "
$ws.Range("D10").Value = "'21.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("D11").Value = "'0.07566"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.98%  "
$ws.Range("D12").Value = "'1.677.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.51%  "
$ws.Range("D13").Value = "'4.463"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.01%  "
$ws.Range("D14").Value = "'0.5616"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.95%  "
$ws.Range("D15").Value = "'66.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").Value = "'0.000008013"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.78%  "
$ws.Range("D17").Value = "'26.039.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.72%  "
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").Value = "'4.817"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.65%  "
$ws.Range("D20").Value = "'187.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("E21").Value = "  -5.14%  "
$ws.Range("D22").Value = "'6.214"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("D23").Value = "'1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("D24").Value = "'149.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").Value = "'0.1255"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.13%  "
$ws.Range("D26").Value = "'7.585"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.02%  "
$ws.Range("D27").Value = "'15.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("D28").Value = "'0.06229"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("D29").Value = "'1.363"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("D30").Value = "'1.285"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.50%  "
$ws.Range("D31").Value = "'3.511"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.95%  "
$ws.Range("D32").Value = "'3.431"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.71%  "
$ws.Range("D33").Value = "'1.635"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.43%  "
$ws.Range("D34").Value = "'1.003"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.65%  "
$ws.Range("D35").Value = "'0.6066"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("D36").Value = "'2.411"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").Value = "'2.753"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").Value = "'6.116"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").Value = "'0.01620"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("D40").Value = "'1.101.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("E41").Value = "  -0.82%  "
$ws.Range("D42").Value = "'1.007"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("D43").Value = "'99.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.34%  "
$ws.Range("D44").Value = "'1.827.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("D45").Value = "'0.00000000110"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.76%  "
$ws.Range("D46").Value = "'56.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.57%  "
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("D48").Value = "'8.045"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.22%  "
$ws.Range("D49").Value = "'0.05232"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.01%  "
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("D51").Value = "'5.980"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.35%  "
